$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.350.99"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").Value = "2.370.91"
$ws.Range("E3").Value = "  +5.11%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.652"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "232.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.78%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.458"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.95%  "

$ws.Range("E10").Value = "  -2.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.25%  "

$ws.Range("D13").Value = "2.718.48"
$ws.Range("E13").Value = "  +4.85%  "

$ws.Range("E14").Value = "  -1.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.94%  "

$ws.Range("D18").Value = "2.366.94"
$ws.Range("E18").Value = "  +4.75%  "

$ws.Range("D19").Value = "43.400.19"
$ws.Range("E19").Value = "  -0.89%  "

$ws.Range("D20").Value = "0.0₃0980"
$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.41%  "

$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +19.68%  "

$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "

$ws.Range("E27").Value = "  +0.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.48%  "

$ws.Range("E31").Value = "  +8.65%  "

$ws.Range("E32").Value = "  -6.88%  "

$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0695"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.71%  "

$ws.Range("E37").Value = "  +10.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.20%  "

$ws.Range("E39").Value = "  -2.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0254"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.21%  "

$ws.Range("E44").Value = "  +8.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.57%  "

$ws.Range("E46").Value = "  +1.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0950"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("E48").Value = "  -1.15%  "

$ws.Range("D49").Value = "1.446.15"
$ws.Range("E49").Value = "  +0.70%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.592.03"
$ws.Range("E50").Value = "  +5.08%  "

$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000205"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.48%  "

